$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.001754667048134761
$ws.Range("C2").Value = 0.0001537489499301437
$ws.Range("D2").Value = 3.082599426703578
$ws.Range("E2").Value = 71517.89157740913
$ws.Range("G2").Value = 71520.97608525184
